$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.860.02"
$ws.Range("E2").Value = "  +2.33%  "
$ws.Range("D3").Value = "2.638.60"
$ws.Range("E3").Value = "  +2.32%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'593.14"
$ws.Range("E5").Value = "  +0.75%  "
$ws.Range("D6").Value = "'155.01"
$ws.Range("E6").Value = "  +3.05%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'0.591"
$ws.Range("E8").Value = "  +0.77%  "
$ws.Range("E9").Value = "  +5.87%  "
$ws.Range("E10").Value = "  +3.46%  "
$ws.Range("E11").Value = "  +0.97%  "
$ws.Range("E12").Value = "  +1.80%  "
$ws.Range("D13").Value = "'28.89"
$ws.Range("E13").Value = "  +5.01%  "
$ws.Range("D14").Value = "'0.0000185"
$ws.Range("E14").Value = "  +19.29%  "
$ws.Range("D15").Value = "3.111.69"
$ws.Range("E15").Value = "  +2.35%  "
$ws.Range("D16").Value = "64.709.26"
$ws.Range("E16").Value = "  +2.38%  "
$ws.Range("D17").Value = "2.669.50"
$ws.Range("E17").Value = "  +3.61%  "
$ws.Range("D18").Value = "'12.53"
$ws.Range("E18").Value = "  +2.90%  "
$ws.Range("D19").Value = "'4.78"
$ws.Range("E19").Value = "  +1.06%  "
$ws.Range("D20").Value = "'350.78"
$ws.Range("E20").Value = "  +1.42%  "
$ws.Range("D21").Value = "'7.24"
$ws.Range("E21").Value = "  +5.78%  "
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").Value = "'67.98"
$ws.Range("E23").Value = "  +1.37%  "
$ws.Range("E24").Value = "  +0.30%  "
$ws.Range("D25").Value = "'9.45"
$ws.Range("E25").Value = "  +3.36%  "
$ws.Range("E26").Value = "  -2.26%  "
$ws.Range("D27").Value = "'8.09"
$ws.Range("E27").Value = "  +0.50%  "
$ws.Range("E28").Value = "  +0.44%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "0.0₃0935"
$ws.Range("E30").Value = "  +9.04%  "
$ws.Range("D31").Value = "'2.09"
$ws.Range("E31").Value = "  +2.80%  "
$ws.Range("D32").Value = "'510.76"
$ws.Range("E32").Value = "  -7.61%  "
$ws.Range("E33").Value = "  +0.70%  "
$ws.Range("D34").Value = "'5.59"
$ws.Range("E34").Value = "  +6.96%  "
$ws.Range("D35").Value = "'6.22"
$ws.Range("E35").Value = "  +3.13%  "
$ws.Range("E36").Value = "  +2.50%  "
$ws.Range("D37").Value = "'164.94"
$ws.Range("E37").Value = "  -1.02%  "
$ws.Range("D38").Value = "'20.09"
$ws.Range("E38").Value = "  +2.88%  "
$ws.Range("E39").Value = "  +4.69%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("D42").Value = "'42.25"
$ws.Range("E42").Value = "  +6.47%  "
$ws.Range("D43").Value = "'164.24"
$ws.Range("E43").Value = "  -0.88%  "
$ws.Range("E44").Value = "  +2.50%  "
$ws.Range("E45").Value = "  +4.35%  "
$ws.Range("D46").Value = "'22.68"
$ws.Range("E46").Value = "  -1.08%  "
$ws.Range("E47").Value = "  +3.44%  "
$ws.Range("D48").Value = "'0.646"
$ws.Range("E48").Value = "  +2.90%  "
$ws.Range("D49").Value = "'0.0253"
$ws.Range("E49").Value = "  +0.79%  "
$ws.Range("D50").Value = "'0.0979"
$ws.Range("E50").Value = "  +1.84%  "
$ws.Range("D51").Value = "'19.27"
$ws.Range("E51").Value = "  +0.80%  "
